$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 635; existing rows 635-685 shift down to 636-686.
$ws.Rows.Item(635).Insert()

# Populate the newly inserted row 635 with the new data record.
$ws.Cells.Item(635, 1).Value = 8
$ws.Cells.Item(635, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(635, 3).Value = "Coquimbo"
$ws.Cells.Item(635, 4).Value = 45194
$ws.Cells.Item(635, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(635, 5).Value = 4
$ws.Cells.Item(635, 6).Value = 100112043
$ws.Cells.Item(635, 7).Value = "Pepino dulce"
$ws.Cells.Item(635, 8).Value = "Sin especificar"
$ws.Cells.Item(635, 9).Value = "Primera"
$ws.Cells.Item(635, 10).Value = 200
$ws.Cells.Item(635, 11).Value = 19000
$ws.Cells.Item(635, 12).Value = 20000
$ws.Cells.Item(635, 13).Value = 19500
$ws.Cells.Item(635, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(635, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(635, 16).Value = 1083
$ws.Cells.Item(635, 17).Value = 18
$ws.Cells.Item(635, 18).Value = "Hortaliza"
